$d = $word.ActiveDocument

# The author re-opened/edited the document and split the erroneously
# concatenated word "afstandgevechten" into the two separate Dutch words
# "afstand gevechten" inside the Mirage paragraph.
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute("afstandgevechten", $true, $false, $false, $false, $false, `
              $true, 1, $false, "afstand gevechten", 2)
